# Update the "want-to-go" headcount (column F) values across all 4 sheets
# to reflect the latest scrape, per commit message:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$updates1 = @{
    3 = 7366
    4 = 3523
    6 = 3857
    8 = 85
    10 = 102
    11 = 156
    12 = 510
    14 = 141
    15 = 373
    19 = 4142
    21 = 413
    22 = 1029
    23 = 537
    24 = 1870
    25 = 116
    26 = 96
    27 = 3053
    28 = 2258
    29 = 64
    32 = 39
    33 = 110
    36 = 4334
    37 = 486
    38 = 324
    41 = 814
    42 = 221
    44 = 1643
    46 = 34
    47 = 610
    48 = 723
}
foreach ($row in $updates1.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates1[$row]
}

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$updates2 = @{
    4 = 441
    12 = 107
    16 = 588
}
foreach ($row in $updates2.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates2[$row]
}

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$updates3 = @{
    2 = 167
    3 = 5
}
foreach ($row in $updates3.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates3[$row]
}

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$updates4 = @{
    2 = 167
    4 = 7366
    5 = 3523
    6 = 3523
    7 = 3857
    8 = 85
    10 = 102
    12 = 156
    13 = 510
    15 = 141
    16 = 373
    20 = 4142
    24 = 413
    25 = 1029
    26 = 537
    27 = 1870
    28 = 116
    29 = 96
    30 = 3053
    31 = 2258
    32 = 64
    35 = 110
    38 = 107
    39 = 4334
    41 = 486
    42 = 324
    44 = 814
    45 = 221
    46 = 1643
    48 = 34
    49 = 610
    50 = 723
}
foreach ($row in $updates4.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates4[$row]
}
